$d = $word.ActiveDocument

# Update the date heading in the first paragraph.
$d.Paragraphs(1).Range.Find.Execute("2024-11-19 Tuesday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-11-20 Wednesday", 2)

# Update the practice-problem table. The table has 20 rows (5 content
# rows with data at rows 1, 5, 9, 13, 17, interleaved with blank spacer
# rows) and 5 columns.
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "13÷8=1, 5"
$t.Cell(1,2).Range.Text = "53÷2=26, 1"
$t.Cell(1,3).Range.Text = "69÷3=23, 0"
$t.Cell(1,4).Range.Text = "51÷8=6, 3"
$t.Cell(1,5).Range.Text = "22÷3=7, 1"

$t.Cell(5,1).Range.Text = "18÷2=9, 0"
$t.Cell(5,2).Range.Text = "16÷4=4, 0"
$t.Cell(5,3).Range.Text = "12÷8=1, 4"
$t.Cell(5,4).Range.Text = "68÷6=11, 2"
$t.Cell(5,5).Range.Text = "92÷9=10, 2"

$t.Cell(9,1).Range.Text = "44÷2=22, 0"
$t.Cell(9,2).Range.Text = "40÷3=13, 1"
$t.Cell(9,3).Range.Text = "61÷8=7, 5"
$t.Cell(9,4).Range.Text = "76÷9=8, 4"
$t.Cell(9,5).Range.Text = "57÷3=19, 0"

$t.Cell(13,1).Range.Text = "18÷3=6, 0"
$t.Cell(13,2).Range.Text = "63÷2=31, 1"
$t.Cell(13,3).Range.Text = "24÷5=4, 4"
$t.Cell(13,4).Range.Text = "65÷5=13, 0"
$t.Cell(13,5).Range.Text = "12÷8=1, 4"

$t.Cell(17,1).Range.Text = "54÷5=10, 4"
$t.Cell(17,2).Range.Text = "55÷2=27, 1"
$t.Cell(17,3).Range.Text = "40÷3=13, 1"
$t.Cell(17,4).Range.Text = "36÷9=4, 0"
$t.Cell(17,5).Range.Text = "64÷5=12, 4"

Write-Output "done"
